$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '57.073.38'
$ws.Cells.Item(2, 5).Value = '  -1.32%  '
$ws.Cells.Item(3, 4).Value = '3.085.30'
$ws.Cells.Item(3, 5).Value = '  -0.23%  '
$ws.Cells.Item(4, 5).Value = '  -0.02%  '
$ws.Cells.Item(5, 4).Value = '520.84'
$ws.Cells.Item(5, 5).Value = '  -0.88%  '
$ws.Cells.Item(6, 4).Value = '135.53'
$ws.Cells.Item(6, 5).Value = '  -3.71%  '
$ws.Cells.Item(7, 5).Value = '  -0.03%  '
$ws.Cells.Item(8, 4).Value = '3.086.46'
$ws.Cells.Item(8, 5).Value = '  -0.16%  '
$ws.Cells.Item(9, 4).Value = '0.453'
$ws.Cells.Item(9, 5).Value = '  +2.40%  '
$ws.Cells.Item(10, 4).Value = '7.35'
$ws.Cells.Item(10, 5).Value = '  +2.75%  '
$ws.Cells.Item(11, 5).Value = '  -1.45%  '
$ws.Cells.Item(12, 4).Value = '0.399'
$ws.Cells.Item(12, 5).Value = '  +1.56%  '
$ws.Cells.Item(13, 2).Value = 'TRON'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(13, 4).Value = '0.136'
$ws.Cells.Item(13, 5).Value = '  +1.78%  '
$ws.Cells.Item(14, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(14, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(14, 4).Value = '3.615.15'
$ws.Cells.Item(14, 5).Value = '  -0.25%  '
$ws.Cells.Item(15, 4).Value = '25.31'
$ws.Cells.Item(15, 5).Value = '  -0.84%  '
$ws.Cells.Item(16, 5).Value = '  -1.89%  '
$ws.Cells.Item(17, 4).Value = '57.179.64'
$ws.Cells.Item(17, 5).Value = '  -1.25%  '
$ws.Cells.Item(18, 4).Value = '3.082.51'
$ws.Cells.Item(18, 5).Value = '  -0.21%  '
$ws.Cells.Item(19, 5).Value = '  -3.30%  '
$ws.Cells.Item(20, 4).Value = '12.49'
$ws.Cells.Item(20, 5).Value = '  -1.61%  '
$ws.Cells.Item(21, 5).Value = '  -1.24%  '
$ws.Cells.Item(22, 4).Value = '347.34'
$ws.Cells.Item(22, 5).Value = '  +1.31%  '
$ws.Cells.Item(23, 2).Value = 'LEO'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Cells.Item(23, 4).Value = '5.81'
$ws.Cells.Item(23, 5).Value = '  +1.63%  '
$ws.Cells.Item(24, 2).Value = 'Dai'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(24, 4).Value = '0.999'
$ws.Cells.Item(24, 5).Value = '  -0.09%  '
$ws.Cells.Item(25, 4).Value = '68.17'
$ws.Cells.Item(25, 5).Value = '  +1.12%  '
$ws.Cells.Item(26, 5).Value = '  -2.52%  '
$ws.Cells.Item(27, 4).Value = '0.166'
$ws.Cells.Item(27, 5).Value = '  -1.83%  '
$ws.Cells.Item(28, 5).Value = '  +0.04%  '
$ws.Cells.Item(29, 4).Value = '0.0₃0863'
$ws.Cells.Item(29, 5).Value = '  -5.63%  '
$ws.Cells.Item(30, 5).Value = '  -0.13%  '
$ws.Cells.Item(31, 4).Value = '7.27'
$ws.Cells.Item(31, 5).Value = '  +0.36%  '
$ws.Cells.Item(32, 5).Value = '  -0.19%  '
$ws.Cells.Item(33, 5).Value = '  -8.36%  '
$ws.Cells.Item(34, 4).Value = '20.83'
$ws.Cells.Item(34, 5).Value = '  -0.64%  '
$ws.Cells.Item(35, 5).Value = '  +6.27%  '
$ws.Cells.Item(36, 4).Value = '159.46'
$ws.Cells.Item(36, 5).Value = '  +0.64%  '
$ws.Cells.Item(37, 4).Value = '1.14'
$ws.Cells.Item(37, 5).Value = '  -3.74%  '
$ws.Cells.Item(38, 4).Value = '6.01'
$ws.Cells.Item(38, 5).Value = '  -2.20%  '
$ws.Cells.Item(39, 4).Value = '25.70'
$ws.Cells.Item(39, 5).Value = '  -1.00%  '
$ws.Cells.Item(40, 5).Value = '  -0.55%  '
$ws.Cells.Item(41, 4).Value = '0.0654'
$ws.Cells.Item(41, 5).Value = '  -1.79%  '
$ws.Cells.Item(42, 4).Value = '1.57'
$ws.Cells.Item(42, 5).Value = '  +2.37%  '
$ws.Cells.Item(43, 4).Value = '4.03'
$ws.Cells.Item(43, 5).Value = '  +0.50%  '
$ws.Cells.Item(44, 4).Value = '0.691'
$ws.Cells.Item(44, 5).Value = '  +1.08%  '
$ws.Cells.Item(45, 4).Value = '2.392.78'
$ws.Cells.Item(45, 5).Value = '  +5.01%  '
$ws.Cells.Item(46, 4).Value = '36.63'
$ws.Cells.Item(46, 5).Value = '  -0.74%  '
$ws.Cells.Item(47, 5).Value = '  -0.02%  '
$ws.Cells.Item(48, 4).Value = '3.124.65'
$ws.Cells.Item(48, 5).Value = '  -0.28%  '
$ws.Cells.Item(49, 5).Value = '  +0.44%  '
$ws.Cells.Item(50, 4).Value = '0.960'
$ws.Cells.Item(50, 5).Value = '  -3.13%  '
$ws.Cells.Item(51, 5).Value = '  -2.41%  '
